$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (RIOT)
$ws.Range("D2").Value = 15.6
$ws.Range("E2").Value = 63.8
$ws.Range("F2").Value = 4.24
$ws.Range("K2").Value = 59.3
$ws.Range("N2").Value = 53.62998959737769

# Row 3 (BTC-USD)
$ws.Range("D3").Value = 92248.23
$ws.Range("E3").Value = 61.7
$ws.Range("F3").Value = 1.54
$ws.Range("K3").Value = 54.1
$ws.Range("N3").Value = 53.62998959737769

# Row 4 (COIN)
$ws.Range("D4").Value = 276.51
$ws.Range("E4").Value = 47.1
$ws.Range("F4").Value = 4.36
$ws.Range("H4").Value = 46
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 51.1
$ws.Range("N4").Value = 53.62998959737769

# Row 5 (MARA)
$ws.Range("D5").Value = 12.42
$ws.Range("E5").Value = 47.3
$ws.Range("F5").Value = 11.83
$ws.Range("K5").Value = 49.3
$ws.Range("N5").Value = 53.62998959737769

# Row 6 (MSTR)
$ws.Range("D6").Value = 188.22
$ws.Range("E6").Value = 40.2
$ws.Range("F6").Value = 7.16
$ws.Range("K6").Value = 38.5
$ws.Range("N6").Value = 53.62998959737769
